$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: cell I2 held a misspelled duplicate of the "cold water" label
# ("Сумма за холодеую воду") instead of the correct text that is already
# used elsewhere in the sheet ("Сумма за холодную воду"). Re-point it at
# the correct text; the now-unused duplicate string is dropped from the
# shared-string table automatically.
$ws.Range("I2").Value = "Сумма за холодную воду"

# Update the active selection to match the current used range.
$ws.Range("A3:XFD17").Select()
